# Replace the Python-style "${:,.2f}".format(X) currency formatting
# expressions used in the Jinja table cells with the Jinja "currency"
# filter: "${:,.2f}".format(X)  ->  X | currency
#
# The document contains exactly two such expressions:
#   {{ "${:,.2f}".format(row.value) }}
#   {{ "${:,.2f}".format(row.paid_by_other_monthly) }}
# which must become:
#   {{ row.value | currency }}
#   {{ row.paid_by_other_monthly | currency }}

$d = $word.ActiveDocument

# ---------- Occurrence 1: row.value ----------

# Remove the leading '"${:,.2f}".format(' wrapper, leaving "row.value)"
$prefix1 = $d.Content.Duplicate
$prefix1.Find.Execute('"${:,.2f}".format(')
$prefix1.Delete()

# Turn the trailing ')' that now follows "row.value" into ' | currency'
$probe1 = $d.Content.Duplicate
$probe1.Find.Execute("row.value)")
$closeParen1 = $d.Range($probe1.End - 1, $probe1.End)
$closeParen1.Text = " | currency"

# ---------- Occurrence 2: row.paid_by_other_monthly ----------

# Remove the leading '"${:,.2f}".format(' wrapper, leaving
# "row.paid_by_other_monthly)"
$prefix2 = $d.Content.Duplicate
$prefix2.Find.Execute('"${:,.2f}".format(')
$prefix2.Delete()

# Turn the trailing ')' that now follows "row.paid_by_other_monthly"
# into ' | currency'
$probe2 = $d.Content.Duplicate
$probe2.Find.Execute("row.paid_by_other_monthly)")
$closeParen2 = $d.Range($probe2.End - 1, $probe2.End)
$closeParen2.Text = " | currency"
